$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -1.8496768167544
$ws.Range("C2").Value = 0.0713815590614557
$ws.Range("B3").Value = -0.201232123211273
$ws.Range("C3").Value = 0.036854147082554

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 0.999415283650382
$ws.Range("C2").Value = 0.0710460548893685
$ws.Range("B3").Value = -0.710160646468684
$ws.Range("C3").Value = 0.0371536156696288

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.37351061111828
$ws.Range("C2").Value = 0.0975451183537697
$ws.Range("B3").Value = 1.17954533219363
$ws.Range("C3").Value = 0.064612489657263

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -1.59135908631409
$ws.Range("C2").Value = 0.0932982144745055
$ws.Range("B3").Value = -0.0889611671755879
$ws.Range("C3").Value = 0.0166836023152949

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.00509532697404409
$ws.Range("B2").Value = -0.000288346995536712
$ws.Range("A3").Value = -0.000288346995536712
$ws.Range("B3").Value = 0.00135822815718253

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.00504754191534316
$ws.Range("B2").Value = -0.00121965787689204
$ws.Range("A3").Value = -0.00121965787689204
$ws.Range("B3").Value = 0.00138039115732648

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00951505011465093
$ws.Range("B2").Value = 0.00140889801365611
$ws.Range("A3").Value = 0.00140889801365611
$ws.Range("B3").Value = 0.00417477381970992

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.00870455682413083
$ws.Range("B2").Value = -0.00105036762195595
$ws.Range("A3").Value = -0.00105036762195595
$ws.Range("B3").Value = 0.000278342586214914
